$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet "Full results": refreshed bootstrap resampling estimates ---
$ws1.Range("H2").Value = 0.564166024789842
$ws1.Range("I2").Value = 0.17988289300563
$ws1.Range("O2").Value = 0.435628374128827
$ws1.Range("F3").Value = 0.576042692642913
$ws1.Range("G3").Value = 0.204450602118496
$ws1.Range("C4").Value = 0.638066291221125
$ws1.Range("D4").Value = 0.361569411877401
$ws1.Range("E4").Value = 0.999635703098526
$ws1.Range("J4").Value = 0.36170117980571
$ws1.Range("K4").Value = 0.204525107058372
$ws1.Range("L4").Value = 0.0118809936245872
$ws1.Range("M4").Value = 0.0739271943231168
$ws1.Range("N4").Value = 0.216406100682959
$ws1.Range("H5").Value = 0.841250028666246
$ws1.Range("I5").Value = 0.0943505841960054
$ws1.Range("O5").Value = 0.158772019642304
$ws1.Range("F6").Value = 0.858224755649418
$ws1.Range("G6").Value = 0.0917999094700253
$ws1.Range("C7").Value = 0.89233621196557
$ws1.Range("D7").Value = 0.107689997855455
$ws1.Range("E7").Value = 1.00002620982103
$ws1.Range("J7").Value = 0.107687175269545
$ws1.Range("K7").Value = 0.0917975033727673
$ws1.Range("L7").Value = 0.0169742821064314
$ws1.Range("M7").Value = 0.0510848443727595
$ws1.Range("N7").Value = 0.108771785479199
$ws1.Range("H8").Value = 0.802766025660298
$ws1.Range("I8").Value = 0.0800516319714449
$ws1.Range("O8").Value = 0.197483910106585
$ws1.Range("F9").Value = 0.845266194618648
$ws1.Range("G9").Value = 0.0827339293747115
$ws1.Range("C10").Value = 0.859701573261984
$ws1.Range("D10").Value = 0.140609863257113
$ws1.Range("E10").Value = 1.0003114365191
$ws1.Range("J10").Value = 0.14056608463557
$ws1.Range("K10").Value = 0.0827081727612653
$ws1.Range("L10").Value = 0.0424869377984927
$ws1.Range("M10").Value = 0.0569178254710148
$ws1.Range("N10").Value = 0.125195110559758
$ws1.Range("H11").Value = 0.766984012917056
$ws1.Range("I11").Value = 0.0666936286405306
$ws1.Range("O11").Value = 0.233285259008839
$ws1.Range("F12").Value = 0.7737235605777
$ws1.Range("G12").Value = 0.0778595744267586
$ws1.Range("C13").Value = 0.827898351653003
$ws1.Range("D13").Value = 0.172452844834897
$ws1.Range("E13").Value = 1.0003511964879
$ws1.Range("J13").Value = 0.172392297681435
$ws1.Range("K13").Value = 0.0778322388826295
$ws1.Range("L13").Value = 0.00673718546760789
$ws1.Range("M13").Value = 0.0608929613274043
$ws1.Range("N13").Value = 0.0845694243502373
$ws1.Range("H14").Value = 0.808729608307432
$ws1.Range("I14").Value = 0.130680500462492
$ws1.Range("O14").Value = 0.192479525848833
$ws1.Range("F15").Value = 0.82900556513295
$ws1.Range("G15").Value = 0.132722287364084
$ws1.Range("C16").Value = 0.828573269432525
$ws1.Range("D16").Value = 0.172924072006996
$ws1.Range("E16").Value = 1.00149734143952
$ws1.Range("J16").Value = 0.17266553142669
$ws1.Range("K16").Value = 0.132523853359574
$ws1.Range("L16").Value = 0.0202456428388175
$ws1.Range("M16").Value = 0.0198139944221428
$ws1.Range("N16").Value = 0.152769496198391

# --- Sheet "For plotting": Index/Outcome columns swapped + refreshed values ---
$ws2.Range("A1").Value = "Index"
$ws2.Range("B1").Value = "Outcome"
$ws2.Range("A2").Value = "Sibcorr"
$ws2.Range("B2").Value = "education"
$ws2.Range("C2").Value = 0.36170117980571
$ws2.Range("D2").Value = 0.322741097211117
$ws2.Range("E2").Value = 0.400661262400303
$ws2.Range("A3").Value = "IOLIB"
$ws2.Range("B3").Value = "education"
$ws2.Range("C3").Value = 0.216406100682959
$ws2.Range("D3").Value = 0.185308575857085
$ws2.Range("E3").Value = 0.247503625508832
$ws2.Range("A4").Value = "IORAD"
$ws2.Range("B4").Value = "education"
$ws2.Range("C4").Value = 0.435628374128827
$ws2.Range("D4").Value = 0.404097397072113
$ws2.Range("E4").Value = 0.467159351185541
$ws2.Range("A5").Value = "Sibcorr"
$ws2.Range("B5").Value = "occupation"
$ws2.Range("C5").Value = 0.172392297681435
$ws2.Range("D5").Value = 0.14204770421408
$ws2.Range("E5").Value = 0.20273689114879
$ws2.Range("A6").Value = "IOLIB"
$ws2.Range("B6").Value = "occupation"
$ws2.Range("C6").Value = 0.0845694243502373
$ws2.Range("D6").Value = 0.0482457849317492
$ws2.Range("E6").Value = 0.120893063768725
$ws2.Range("A7").Value = "IORAD"
$ws2.Range("B7").Value = "occupation"
$ws2.Range("C7").Value = 0.233285259008839
$ws2.Range("D7").Value = 0.204682421394656
$ws2.Range("E7").Value = 0.261888096623022
$ws2.Range("A8").Value = "Sibcorr"
$ws2.Range("B8").Value = "income"
$ws2.Range("C8").Value = 0.14056608463557
$ws2.Range("D8").Value = 0.0512911993422382
$ws2.Range("E8").Value = 0.229840969928901
$ws2.Range("A9").Value = "IOLIB"
$ws2.Range("B9").Value = "income"
$ws2.Range("C9").Value = 0.125195110559758
$ws2.Range("D9").Value = 0.0360863304357876
$ws2.Range("E9").Value = 0.214303890683728
$ws2.Range("A10").Value = "IORAD"
$ws2.Range("B10").Value = "income"
$ws2.Range("C10").Value = 0.197483910106585
$ws2.Range("D10").Value = 0.114939916178892
$ws2.Range("E10").Value = 0.280027904034278
$ws2.Range("A11").Value = "Sibcorr"
$ws2.Range("B11").Value = "wealth"
$ws2.Range("C11").Value = 0.17266553142669
$ws2.Range("D11").Value = 0.0560058845213903
$ws2.Range("E11").Value = 0.28932517833199
$ws2.Range("A12").Value = "IOLIB"
$ws2.Range("B12").Value = "wealth"
$ws2.Range("C12").Value = 0.152769496198391
$ws2.Range("D12").Value = 0.0390740685062593
$ws2.Range("E12").Value = 0.266464923890523
$ws2.Range("A13").Value = "IORAD"
$ws2.Range("B13").Value = "wealth"
$ws2.Range("C13").Value = 0.192479525848833
$ws2.Range("D13").Value = 0.0836957137848644
$ws2.Range("E13").Value = 0.301263337912802
$ws2.Range("A14").Value = "Sibcorr"
$ws2.Range("B14").Value = "health_pc"
$ws2.Range("C14").Value = 0.107687175269545
$ws2.Range("D14").Value = 0.0522538340451524
$ws2.Range("E14").Value = 0.163120516493937
$ws2.Range("A15").Value = "IOLIB"
$ws2.Range("B15").Value = "health_pc"
$ws2.Range("C15").Value = 0.108771785479199
$ws2.Range("D15").Value = 0.0545047685720089
$ws2.Range("E15").Value = 0.163038802386389
$ws2.Range("A16").Value = "IORAD"
$ws2.Range("B16").Value = "health_pc"
$ws2.Range("C16").Value = 0.158772019642304
$ws2.Range("D16").Value = 0.104751073568448
$ws2.Range("E16").Value = 0.21279296571616
